# Commit: "Fruta / hortaliza, semanal"
# A new weekly price observation is inserted as a new data row right after
# the header/first rows, at row 19. Every existing row from 19 downward
# shifts down by one (old row 19 -> new row 20, ..., old row 126 -> new row 127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 - this pushes rows 19..126 down to 20..127
# and copies formatting (incl. the date number-format) from the row above,
# matching the canonical OOXML (dimension grows from R126 to R127).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new observation.
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44473
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112006
$ws.Range("G19").Value = "Repollo"
$ws.Range("H19").Value = "Crespo record"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 550
$ws.Range("M19").Value = 525
$ws.Range("N19").Value = "`$/unidad"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 525
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"
